$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tables")
$lo = $ws.ListObjects.Item(1)

# Insert a new row for the "setup" table entry above the current row 3
# (area_users), shifting all existing rows down by one.
$ws.Rows.Item(3).Insert()

# Grow the table/list object so it covers the newly inserted row as well
# as the extra row appended at the bottom (A1:E23 -> A1:E24).
$lo.Resize($ws.Range("A1:E24"))

# Populate the new row with the "setup" table definition.
$ws.Range("A3").Value = "setup"
$ws.Range("B3").Value = "name,value"
$ws.Range("C3").Value = "APP"
$ws.Range("D3").Value = 24

# Update the "users" row (now row 24) with the expanded list of fields.
$ws.Range("B24").Value = "name,code,phone,outstanding_normal,outstanding_overdue,outstanding_critical"

# Reflect the user's final selection/scroll position on the sheet.
$ws.Range("B25").Select()
